# Applies the commit's two textual simplifications inside the header table:
#   1. "Fecha del doc:" -> "Fecha:"          (collapses 4 runs into one)
#   2. "{NumDoc}" (split across 3 runs)       (collapses into a single run,
#      "{" + "NumDoc" + "}")                  dropping the spell-check markers)
#
# Word's Find/Replace (wdReplaceAll) merges the matched runs into a single
# run carrying the formatting of the first matched run, and drops any
# w:proofErr spell-check bookmarks that fell inside the replaced range -
# exactly what the diff shows.

$d = $word.ActiveDocument

# 1) "Fecha del doc:" -> "Fecha:"
$d.Content.Find.Execute(
    "Fecha del doc:",  # FindText
    $true,              # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,               # Forward
    1,                   # Wrap (wdFindContinue)
    $false,              # Format
    "Fecha:",            # ReplaceWith
    2                    # Replace (wdReplaceAll)
) | Out-Null

# 2) "{NumDoc}" -> "{NumDoc}" (merges the 3 split runs into one)
$d.Content.Find.Execute(
    "{NumDoc}",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "{NumDoc}",
    2
) | Out-Null
